$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A136").Value = "2023-12-09 10:18:19"
$ws.Range("B136").Value = 0.0006000000000000001

$ws.Range("A137").Value = "2023-12-09 10:18:24"
$ws.Range("B137").Value = 0.0002
